$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D; this shifts the existing D/E columns to E/F
# and auto-adjusts the formula in the (now) E3 cell that referenced D2.
$ws.Columns("D").Insert()

# New header for the inserted column.
$ws.Range("D1").Value = "errors"

# D2: a formula-driven error (#DIV/0!), per BoolErr handling of formula errors.
$ws.Range("D2").Formula = "=1/0"

# D3: a literal (non-formula) error value (#N/A).
$ws.Range("D3").Value = "#N/A"

# D4: another formula-driven error (#VALUE!), adding a number to a string.
$ws.Range("D4").Formula = "=F2+F4"

# Flip the existing boolean column (C) values.
[void]($ws.Range("C2").Value = $true)
[void]($ws.Range("C3").Value = $false)

# Update the selected cell to match the new layout.
[void]$ws.Range("E3").Select()
